$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 4
$ws.Range("X5").Value = 8
$ws.Range("Z5").Value = 17
$ws.Range("AA5").Value = 21
$ws.Range("AJ5").Value = 41
$ws.Range("AO5").Value = 12
$ws.Range("AP5").Value = 29
$ws.Range("AR5").Value = 81
